# Update Handback status timestamps ("Generate Report for Handback")
$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Row for 099ae328-bff9-4406-bf11-3fcf9e46e4bf, column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-10-19 23:36:01"

# --- zh-cn sheet ---
# Row for 099ae328-bff9-4406-bf11-3fcf9e46e4bf.1e31b8ac3255d01ec9ea086a387e06fc50a9d65a.zh-cn.xlf
# H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-10-19 23:35:49"
$wsZhCn.Range("K3").Value = "2016-10-19 23:36:35"

# --- de-de sheet ---
# Row for 099ae328-bff9-4406-bf11-3fcf9e46e4bf.1e31b8ac3255d01ec9ea086a387e06fc50a9d65a.de-de.xlf
# K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = "2016-10-19 23:36:53"
